$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New translation rows: Key (A) / es_ES (B) / en_GB (C)
$rows = @(
    @("mail_ca_enabled_body", "{0}, Tu cuenta como CA ha sido habilitada exitosamente, ahora puedes dar de alta nuevos cursos y emitir nuevos certificados", "{0}, Your account as CA has been successfully enabled, now you can register new courses and issuer new certificates"),
    @("mail_ca_disabled_body", "{0}, Tu cuenta como CA ha sido deshabilitada, no podrás registrar nuevos cursos ni emitir nuevos certificados hasta que sea activada de nuevo.", "{0}, Your CA account has been disabled, you will not be able to register new courses or issue new certificates until it is activated again."),
    @("mail_certificate_disabled_body", "{0}, Tú certificado {1} ha sido deshabilitado, este no podrá ser renovado ni utilizado hasta que sea habilitado de nuevo.", "{0}, Your certificate {1} has been disabled, it cannot be renewed or used until it is enabled again."),
    @("mail_certificate_enabled_body", "{0}, Tú certificado {1} ha sido habilitado, podrás renovar y utilizar tu certificado con normalidad.", "{0}, Your certificate {1} has been enabled, you can renew and use your certificate normally."),
    @("mail_certificate_renewed_body", "{0}, Tú certificado {1} ha sido renovado, podrás seguir utilizándolo con normalidad.", "{0}, Your certificate {1} has been renewed, you can continue to use it normally."),
    @("mail_certificate_request_accepted_body", "{0}, Tú solicitud de emisión del certificado {1} ha sido aceptada, en breve tú certificado será generado y almacenado en TCS.", "{0}, Your request to issue the certificate {1} has been accepted, shortly your certificate will be generated and stored in TCS."),
    @("mail_certificate_request_rejected_body", "{0}, Tú solicitud de emisión del certificado {1} ha sido rechazada.", "{0}, Your request to issue the certificate {1} has been rejected."),
    @("mail_certificate_changed_to_visible_body", "{0}, La visibilidad de tú certificado {1} ha sido actualizada, tú certificado ahora es visible.", "{0}, The visibility of your certificate {1} has been updated, your certificate is now visible."),
    @("mail_certificate_changed_to_invisible_body", "{0}, La visibilidad de tú certificado {1} ha sido actualizada, tú certificado ahora no es visible.", "{0}, The visibility of your certificate {1} has been updated, your certificate is now not visible.")
)

$startRow = 66
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

$ws.Range("C74").Select()
